$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C41").Value = 126945
$ws.Range("E41").Value = 662707757

$ws.Range("C48").Value = 150636
$ws.Range("E48").Value = 275744616

$ws.Range("C67").Value = 27105
$ws.Range("E67").Value = 168714702

$ws.Range("C69").Value = 17891
$ws.Range("E69").Value = 103836984

$ws.Range("C72").Value = 331325
$ws.Range("E72").Value = 635436311

$ws.Range("C81").Value = 88356
$ws.Range("E81").Value = 499673318

$ws.Range("C104").Value = 319351
$ws.Range("D104").Value = 67905
$ws.Range("E104").Value = 561273145

$ws.Range("C121").Value = 1306327
$ws.Range("D121").Value = 220386
$ws.Range("E121").Value = 2275257402

$ws.Range("C129").Value = 633686
$ws.Range("E129").Value = 3433172989

$ws.Range("C132").Value = 585952
$ws.Range("E132").Value = 3470848013

$ws.Range("C186").Value = 236835
$ws.Range("E186").Value = 1189990205

$ws.Range("C189").Value = 100471
$ws.Range("E189").Value = 556436057

$ws.Range("C204").Value = 265654
$ws.Range("E204").Value = 1271493760

$ws.Range("C237").Value = 283322
$ws.Range("E237").Value = 1438427690

$ws.Range("C240").Value = 205926
$ws.Range("E240").Value = 1069823635
